$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo: "Nun\nez" -> "Nu\~{n}ez" for author A. Nunez (row 77)
$ws.Range("A77").Value = "Nu\~{n}ez"

# Delete rows for removed authors (delete bottom-to-top to keep row numbers stable)
# Row 110: Yahlali, N.
$ws.Rows.Item(110).Delete()
# Row 88: Ripoll, L.
$ws.Rows.Item(88).Delete()
# Row 84: Perez, J.
$ws.Rows.Item(84).Delete()
# Row 63: Losada, M.
$ws.Rows.Item(63).Delete()
# Row 45: Gutierrez, R.M.
$ws.Rows.Item(45).Delete()
# Row 27: Diaz, J.
$ws.Rows.Item(27).Delete()

$ws.Range("A105").Select()
